# Generate Report for Handback
#
# - Marks both source files as handed back (Status column -> "Handed back:
#   in sync with en-US") on the Overview sheet as well as on each language
#   sheet (zh-cn / de-de).
# - Fills in the "Latest Target File" (F) and "Latest Handback File" (G)
#   columns on the zh-cn / de-de sheets with hyperlinked file names, now
#   that a handback package exists.
# - Stamps the "Latest Handback DateTime" (H) column with the actual
#   handback timestamp (previously the zero-date placeholder).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: Status columns for both language columns / both files.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusText
$overview.Range("C2").Value = $statusText
$overview.Range("B3").Value = $statusText
$overview.Range("C3").Value = $statusText

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status column
$zhcn.Range("C2").Value = $statusText
$zhcn.Range("C3").Value = $statusText

# Latest Target File / Latest Handback File hyperlinks for row 2
# (source file d3a31d31-5aa2-4599-8359-ef1e57e533f9.md)
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/5d69e0515317525ad6b4abf21c2ad4349346bb4f/e2e/d3a31d31-5aa2-4599-8359-ef1e57e533f9.md", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbae56b213e639f609cc5aec3c1ce642490ad52a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.zh-cn.xlf", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.zh-cn.xlf")

# Latest Target File / Latest Handback File hyperlinks for row 3
# (source file d6944fb7-3bb4-420d-9661-1aa6ff230e90.md)
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/5d69e0515317525ad6b4abf21c2ad4349346bb4f/e2e/d6944fb7-3bb4-420d-9661-1aa6ff230e90.md", "", "", "d6944fb7-3bb4-420d-9661-1aa6ff230e90.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fbae56b213e639f609cc5aec3c1ce642490ad52a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d6944fb7-3bb4-420d-9661-1aa6ff230e90.f0d9538132d0002b3f380952617c7b65f5b4428a.zh-cn.xlf", "", "", "d6944fb7-3bb4-420d-9661-1aa6ff230e90.f0d9538132d0002b3f380952617c7b65f5b4428a.zh-cn.xlf")

# Latest Handback DateTime
$zhcn.Range("H2").Value = "2016-03-11 20:33:17"
$zhcn.Range("H3").Value = "2016-03-11 20:33:17"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

# Status column
$dede.Range("C2").Value = $statusText
$dede.Range("C3").Value = $statusText

# Latest Target File / Latest Handback File hyperlinks for row 2
# (source file d3a31d31-5aa2-4599-8359-ef1e57e533f9.md)
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/5d69e0515317525ad6b4abf21c2ad4349346bb4f/e2e/d3a31d31-5aa2-4599-8359-ef1e57e533f9.md", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d13600a736fab76aae76a58577640f09ad8d5b27/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.de-de.xlf", "", "", "d3a31d31-5aa2-4599-8359-ef1e57e533f9.4b8fd5364e0172a0c45ee5a376e1f8497edbe73a.de-de.xlf")

# Latest Target File / Latest Handback File hyperlinks for row 3
# (source file d6944fb7-3bb4-420d-9661-1aa6ff230e90.md)
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/5d69e0515317525ad6b4abf21c2ad4349346bb4f/e2e/d6944fb7-3bb4-420d-9661-1aa6ff230e90.md", "", "", "d6944fb7-3bb4-420d-9661-1aa6ff230e90.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d13600a736fab76aae76a58577640f09ad8d5b27/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d6944fb7-3bb4-420d-9661-1aa6ff230e90.f0d9538132d0002b3f380952617c7b65f5b4428a.de-de.xlf", "", "", "d6944fb7-3bb4-420d-9661-1aa6ff230e90.f0d9538132d0002b3f380952617c7b65f5b4428a.de-de.xlf")

# Latest Handback DateTime (distinct timestamp from zh-cn's handback)
$dede.Range("H2").Value = "2016-03-11 20:33:22"
$dede.Range("H3").Value = "2016-03-11 20:33:22"
